# Update countries & provincias Spain
# - Refresh statistics for several countries (value-only updates).
# - The table is kept sorted by "Casos totales" (col B) descending; since some
#   countries' totals changed rank, their row positions are re-used but the
#   country name + stats in those rows are updated to reflect the new order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 116050
$ws.Range("C4").Value = 11924
$ws.Range("D4").Value = 3224
$ws.Range("E4").Value = 110889
$ws.Range("G4").Value = 241
$ws.Range("H4").Value = 1937

$ws.Range("A18").Value = "Canada"
$ws.Range("B18").Value = 5434
$ws.Range("C18").Value = 677
$ws.Range("D18").Value = 354
$ws.Range("E18").Value = 5025
$ws.Range("F18").Value = 120
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 55

$ws.Range("A19").Value = "Portugal"
$ws.Range("B19").Value = 5170
$ws.Range("C19").Value = 902
$ws.Range("D19").Value = 43
$ws.Range("E19").Value = 5027
$ws.Range("F19").Value = 89
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 100

$ws.Range("B20").Value = 3977
$ws.Range("C20").Value = 206
$ws.Range("E20").Value = 3950

$ws.Range("E25").Value = 2519
$ws.Range("F25").Value = 250
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 11

$ws.Range("E39").Value = 1137
$ws.Range("H39").Value = 2

$ws.Range("A58").Value = "Egipto"
$ws.Range("B58").Value = 576
$ws.Range("C58").Value = 40
$ws.Range("D58").Value = 121
$ws.Range("E58").Value = 419
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 6
$ws.Range("H58").Value = 36

$ws.Range("A59").Value = "Hong Kong"
$ws.Range("B59").Value = 560
$ws.Range("C59").Value = 41
$ws.Range("D59").Value = 112
$ws.Range("E59").Value = 444
$ws.Range("F59").Value = 5
$ws.Range("H59").Value = 4

$ws.Range("A60").Value = "Colombia"
$ws.Range("B60").Value = 539
$ws.Range("D60").Value = 10
$ws.Range("E60").Value = 523
$ws.Range("H60").Value = 6

$ws.Range("B71").Value = 331
$ws.Range("C71").Value = 38
$ws.Range("E71").Value = 314

$ws.Range("A188").Value = "Libia"
$ws.Range("C188").Value = 2

$ws.Range("A190").Value = "Republica de Africa Central"

$ws.Range("A192").Value = "Butan"

$ws.Range("A193").Value = "Liberia"

$ws.Range("A194").Value = "Somalia"
$ws.Range("E194").Value = 3
$ws.Range("H194").Value = 0

$ws.Range("A196").Value = "Gambia"
$ws.Range("B196").Value = 3
$ws.Range("H196").Value = 1

$ws.Range("A197").Value = "Islas Turcas y Caicos"

$ws.Range("A198").Value = "Islas Virgenes Britanicas"

$ws.Range("A199").Value = "Belice"

$ws.Range("A200").Value = "Anguila"

$ws.Range("A201").Value = "Guinea-Bisau"

$ws.Range("A202").Value = "San Cristobal y Nieves"
$ws.Range("B202").Value = 2
$ws.Range("E202").Value = 2

$ws.Range("A203").Value = "Papua Nueva Guinea"

$ws.Range("A204").Value = "Timor Oriental"

$ws.Range("A205").Value = "San Vicente y las Granadinas"
